# Weekly update: insert a new reporting week (date 44595) of 3 rows
# ("1a amarillo", "2a amarillo", "3a amarillo") at the top of the
# "Comercializadora del Agro de Limarí - Limón" block (old row 483),
# pushing the existing rows 483-496 down to 486-499.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three fresh blank rows above the current row 483. Doing this three
# times at the same index shifts everything else down by three rows total
# while the newly inserted rows inherit the formatting (incl. the date
# number format in column D) from the row that used to be there.
$ws.Rows.Item(483).Insert()
$ws.Rows.Item(483).Insert()
$ws.Rows.Item(483).Insert()

# Data for the three new rows (row 483..485), columns A..T.
$newRows = @(
    @{ Row = 483; L = "1a amarillo"; M = 750; N = 10800; O = 11000; P = 10900; S = 681 },
    @{ Row = 484; L = "2a amarillo"; M = 600; N = 8800;  O = 9000;  P = 8900;  S = 556 },
    @{ Row = 485; L = "3a amarillo"; M = 450; N = 5800;  O = 6000;  P = 5900;  S = 369 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44595
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102003
    $ws.Cells.Item($row, 10).Value = "Limón"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/malla 16 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 16
}
